# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# - Update VALOR MORA total (E11) and Cant. Periodos (F13)
# - Insert two new "Estado de Cuenta" detail rows (period 2509) for the two
#   existing employees, pushing the closing/signature block down
# - Keep the previously-last detail row (period 2508 / JENADIS) as a normal
#   interior row and move the "last row" (bottom border) styling to the new
#   final detail row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -------------------------------------------------
$ws.Range("E11").Value = 596587
$ws.Range("F13").Value = 7

# --- Make room for two new detail rows after row 25 (current last row) -----
$ws.Rows("26:27").Insert()

# Row 26 should look like a normal interior detail row (copy format from 24)
$ws.Range("B24:J24").Copy() | Out-Null
$ws.Range("B26:J26").PasteSpecial(-4122) | Out-Null

# Row 27 becomes the new final row, so it gets the "bottom border" look that
# used to belong to the old row 25
$ws.Range("B25:J25").Copy() | Out-Null
$ws.Range("B27:J27").PasteSpecial(-4122) | Out-Null

# Row 25 is no longer the last row, so restyle it like a normal interior row
$ws.Range("B24:J24").Copy() | Out-Null
$ws.Range("B25:J25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New row 26: SINDY PAOLA IRIARTE CASTILLO, periodo 2509 -----------------
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1044927134"
$ws.Range("D26").Value = "SINDY PAOLA IRIARTE CASTILLO"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2509"
$ws.Range("F26").Value = 56940
$ws.Range("G26").Value = 1423500

# --- New row 27: JENADIS DE LA ROSA MARQUEZ, periodo 2509 -------------------
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1063148285"
$ws.Range("D27").Value = "JENADIS DE LA ROSA MARQUEZ"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2509"
$ws.Range("F27").Value = 56940
$ws.Range("G27").Value = 1423500
